$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.311.44"
$ws.Range("E2").Value = "  +2.58%  "
$ws.Range("D3").Value = "2.525.00"
$ws.Range("E3").Value = "  +1.92%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.61"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.45"
$ws.Range("E6").Value = "  +1.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.532"
$ws.Range("E7").Value = "  +2.31%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.553"
$ws.Range("E9").Value = "  +3.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.46"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.60"
$ws.Range("E11").Value = "  +13.81%  "
$ws.Range("E12").Value = "  +2.21%  "
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("D15").Value = "2.919.30"
$ws.Range("E15").Value = "  +2.51%  "
$ws.Range("D16").Value = "2.533.36"
$ws.Range("E16").Value = "  +2.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.855"
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("D18").Value = "48.142.29"
$ws.Range("E18").Value = "  +2.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.44"
$ws.Range("E19").Value = "  +6.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.65"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("E21").Value = "  +2.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.71"
$ws.Range("E22").Value = "  -0.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.94"
$ws.Range("E23").Value = "  +2.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.44"
$ws.Range("E24").Value = "  +8.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.57"
$ws.Range("E25").Value = "  +1.03%  "
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.09"
$ws.Range("E27").Value = "  +2.07%  "
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("E29").Value = "  -2.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.144"
$ws.Range("E30").Value = "  +5.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.82"
$ws.Range("E31").Value = "  +2.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.66"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.78"
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.40"
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0788"
$ws.Range("E36").Value = "  +1.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.99"
$ws.Range("E37").Value = "  +2.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.72"
$ws.Range("E38").Value = "  +1.94%  "
$ws.Range("E39").Value = "  +2.75%  "
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "121.05"
$ws.Range("E41").Value = "  +1.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.11"
$ws.Range("E42").Value = "  +0.79%  "
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("E44").Value = "  +2.64%  "
$ws.Range("D45").Value = "2.017.21"
$ws.Range("E45").Value = "  +2.05%  "
$ws.Range("E46").Value = "  +5.90%  "
$ws.Range("E47").Value = "  +8.83%  "
$ws.Range("E48").Value = "  +2.01%  "
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("E50").Value = "  +2.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.22"
$ws.Range("E51").Value = "  +3.34%  "
